$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Highlight (green) the intro text run in paragraph 1.
# ---------------------------------------------------------------------
$introRng = $d.Range(0, 0)
$introRng.Find.Execute("Cambiar el date del formulario de registro: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$introRng.HighlightColorIndex = 4

# ---------------------------------------------------------------------
# 2) Highlight (green) + color=auto the hyperlink run. Directly setting
#    HighlightColorIndex on a hyperlink-wrapped run is a no-op in this
#    engine, so drive it through Find/Replace with formatting instead
#    (re-asserting the Hipervnculo character style so it isn't lost).
# ---------------------------------------------------------------------
$hlRng = $d.Range(0, 0)
$hlRng.Find.ClearFormatting()
$hlRng.Find.Replacement.ClearFormatting()
$hlRng.Find.Replacement.Font.HighlightColorIndex = 4
$hlRng.Find.Replacement.Font.Color = -16777216
$hlRng.Find.Replacement.Style = "Hipervnculo"
$hlRng.Find.Execute("http://www.eyecon.ro/datepicker/", $true, $false, $false, $false, $false, $true, 1, $false, "http://www.eyecon.ro/datepicker/", 2, $true, $false, 0, $false)

# ---------------------------------------------------------------------
# 3) Insert the new paragraphs (with blank-line separators) right
#    before the bookmark paragraph ("_GoBack"), keeping it and the
#    section properties where they already are.
# ---------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item(2)

$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()

$d.Paragraphs.Item(3).Range.Text = "Falta poner en registro controler en base a las exepciones de la bd una notificación hacia registro.php que muestre cuando un usuario se ha repetido"
$d.Paragraphs.Item(5).Range.Text = "En los formularios de creación de elementos falta un botón de comprobar disponibilidad"

# ---------------------------------------------------------------------
# 4) Drop the trailing empty paragraph that used to sit right after the
#    bookmark paragraph (the document now ends with the bookmark
#    paragraph immediately followed by the section break). Deleting a
#    zero-length Range is a no-op, so the deleted range has to reach
#    back and swallow the previous paragraph mark too.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$prevEnd = $d.Paragraphs.Item($n - 1).Range.End
$lastEnd = $d.Paragraphs.Item($n).Range.End
$trailRng = $d.Range($prevEnd - 1, $lastEnd)
$trailRng.Delete()
